# "removed mixed from aws"
# The AWS sheet tracked an extra "Mixed" sentiment bucket (column D) that the
# other sheets (Watson/Google/Microsoft/Google15) don't have. This edit drops
# that bucket: its label/count are cleared, the 86 "Mixed" comments are folded
# back into the Neutral row's True/False-Positive counts, the Mismatches total
# is adjusted accordingly, and the per-class Precision/Recall/F-score cells
# (previously hard-coded numbers) are rewritten as live formulas - matching
# the pattern already used on the other per-tool sheets.

$wb = $excel.ActiveWorkbook
$aws = $wb.Worksheets.Item("AWS")

# --- Remove the "Mixed" column from the summary block (row 1/2) ---
$aws.Range("D1").ClearContents()
$aws.Range("D2").ClearContents()

# Mismatches total no longer includes the old 86 "Mixed" comments, but does
# include the extra comments that moved into Neutral's False Positives.
$aws.Range("G2").Value = 1906

# Accuracy figure tied to the old Mixed-inclusive formula is cleared.
$aws.Range("I2").ClearContents()

# --- Neutral class True Positives / False Positives absorb the old Mixed bucket ---
$aws.Range("B15").Value = 2396
$aws.Range("C15").Value = 271

# --- Replace the hard-coded Precision / Recall / F-score cells with formulas ---
$aws.Range("D13").Formula = "=B13/(B13+C13)"
$aws.Range("E13").Formula = "=B13/A6"
$aws.Range("F13").Formula = "=2*D13*E13/(D13+E13)"

$aws.Range("D14").Formula = "=B14/(B14+C14)"
$aws.Range("E14").Formula = "=B14/B6"
$aws.Range("F14").Formula = "=2*D14*E14/(D14+E14)"

$aws.Range("D15").Formula = "=B15/(B15+C15)"
$aws.Range("E15").Formula = "=B15/C6"
$aws.Range("F15").Formula = "=2*D15*E15/(D15+E15)"

# These formula cells pick up the default (11pt) font instead of the sheet's
# 12pt style, same as the already-live formula cells on the sibling sheets.
foreach ($addr in @("D13","E13","F13","D14","E14","F14","D15","E15","F15")) {
    $aws.Range($addr).Font.Size = 11
}

# --- Update the AWS bar chart so it no longer plots the removed "Mixed" point ---
$chartObj = $aws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(,AWS!`$A`$1:`$C`$1,AWS!`$A`$2:`$C`$2,1)"

# --- The AWS sheet was the active/selected sheet when this edit was made ---
$aws.Activate()
$aws.Range("G13").Select()
